$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- 1) Add two new table columns: cross_section, cross_section_units ---
$colCS  = $lo.ListColumns.Add()
$colCSU = $lo.ListColumns.Add()

$ws.Range("S1").Value = "cross_section"
$ws.Range("T1").Value = "cross_section_units"

# Header formatting to match the rest of the header row (bold, white, centered, bordered)
$hdr = $ws.Range("S1:T1")
$hdr.Font.Bold = $true
$hdr.Font.ThemeColor = 2
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1
$hdr.Borders.Item(8).LineStyle = -4142

# --- 2) Backfill existing data rows (2-15) with default cross_section values ---
for ($r = 2; $r -le 15; $r++) {
  $ws.Cells.Item($r, 19).Value = 0
  $ws.Cells.Item($r, 20).Value = "mm"
}

# --- 3) Add three new table rows for the Almelec conductors ---
$row16 = $lo.ListRows.Add()
$row17 = $lo.ListRows.Add()
$row18 = $lo.ListRows.Add()

# Row 16: Almelec_50
$ws.Range("A16").Value = "Almelec_50"
$ws.Range("B16").Value = 0.064100000000000004
$ws.Range("C16").Value = 0.315
$ws.Range("D16").Value = 0.31
$ws.Range("E16").Value = 1.5
$ws.Range("F16").Value = 16
$ws.Range("G16").Value = 0.94
$ws.Range("H16").Value = 0.048675000000000003
$ws.Range("I16").Value = "ft"
$ws.Range("J16").Value = 14.872199999999999
$ws.Range("K16").Value = 0.82399999999999995
$ws.Range("L16").Value = 205
$ws.Range("M16").Value = "km"
$ws.Range("N16").Value = 0.0020799999999999998
$ws.Range("O16").Value = "in"
$ws.Range("P16").Value = 30
$ws.Range("Q16").Value = "1/3 Neutral"
$ws.Range("R16").Value = "ALMELEC"
$ws.Range("S16").Value = 50
$ws.Range("T16").Value = "mm"

# Row 17: Almelec_120
$ws.Range("A17").Value = "Almelec_120"
$ws.Range("B17").Value = 0.064100000000000004
$ws.Range("C17").Value = 0.315
$ws.Range("D17").Value = 0.49
$ws.Range("E17").Value = 1.69
$ws.Range("F17").Value = 16
$ws.Range("G17").Value = 1.1200000000000001
$ws.Range("H17").Value = 0.054840500000000007
$ws.Range("I17").Value = "ft"
$ws.Range("J17").Value = 14.872199999999999
$ws.Range("K17").Value = 0.32700000000000001
$ws.Range("L17").Value = 295
$ws.Range("M17").Value = "km"
$ws.Range("N17").Value = 0.0020799999999999998
$ws.Range("O17").Value = "in"
$ws.Range("P17").Value = 30
$ws.Range("Q17").Value = "1/3 Neutral"
$ws.Range("R17").Value = "ALMELEC"
$ws.Range("S17").Value = 120
$ws.Range("T17").Value = "mm"

# Row 18: Almelec_240
$ws.Range("A18").Value = "Almelec_240"
$ws.Range("B18").Value = 0.064100000000000004
$ws.Range("C18").Value = 0.315
$ws.Range("D18").Value = 0.69
$ws.Range("E18").Value = 1.89
$ws.Range("F18").Value = 25
$ws.Range("G18").Value = 1.3199999999999998
$ws.Range("H18").Value = 0.061330500000000003
$ws.Range("I18").Value = "ft"
$ws.Range("J18").Value = 14.872199999999999
$ws.Range("K18").Value = 0.16400000000000001
$ws.Range("L18").Value = 390
$ws.Range("M18").Value = "km"
$ws.Range("N18").Value = 0.0020799999999999998
$ws.Range("O18").Value = "in"
$ws.Range("P18").Value = 30
$ws.Range("Q18").Value = "1/3 Neutral"
$ws.Range("R18").Value = "ALMELEC"
$ws.Range("S18").Value = 240
$ws.Range("T18").Value = "mm"

# --- 4) Highlight the newly-sourced strand/resistance figures in red ---
$ws.Range("B16:B18").Font.Color = 255
$ws.Range("J16:J18").Font.Color = 255

# --- 5) Keep the same "General"-style numeric formatting used elsewhere in F/L columns ---
$ws.Range("F17").NumberFormat = "General"
$ws.Range("F18").NumberFormat = "General"
$ws.Range("L17").NumberFormat = "General"
$ws.Range("L18").NumberFormat = "General"

# --- 6) Match the final selection left behind in the workbook ---
$ws.Range("B17").Select()
